# Re-fit of the (d,p) spectroscopic-factor table: energies (B), angular
# momenta (C), spectroscopic factors (D) and errors (E) are replaced with
# the "moved back to 116Cd (d,p)" values, and the table shrinks from 51
# levels (rows 2-53) down to 47 levels (rows 2-48).
$data = @(
    @(1222,0,0.01654890364994536,0.002249365544652767),
    @(1237,2,0.0005321444568159682,0.0001900515917199887),
    @(1257,0,0.01203606009123783,0.00174206132899495),
    @(1279,2,0.01442844011991008,0.0006035425932511405),
    @(1316.961967096084,0,0.004580565564309971,0.001263604293602751),
    @(1343,2,0.004212753520500596,0.0004100910506681996),
    @(1355,2,0.01352227287378613,0.0005952032076494585),
    @(1423.282918228582,1,0.002279683569842992,0.0004844350193860269),
    @(1440.359726841536,1,0.002539610150580014,0.0005118901462955945),
    @(1476,2,0.01843077524525283,0.001182626249695394),
    @(1484.351379665093,5,0.03066038090138969,0.002647382286666178),
    @(1648.526094901215,2,0.0010229163193118,0.0002292743474319551),
    @(1672.06961683155,2,0.003792593738153586,0.0005443074346424128),
    @(1682.013801340207,1,0.005524460759057719,0.0005288007084618736),
    @(1692.004774397466,2,0.003308365968041293,0.001750458184148833),
    @(1743.34190550635,2,0.001213154131293051,0.0002772923728669831),
    @(1844.463090411462,1,0.01324288828617537,0),
    @(1844.463090411462,2,0.01653707255420836,0),
    @(1880.855063637229,2,0.0180723617287715,0.0004512474709394471),
    @(1943,0,[double]"1.399646823909309e-20",0),
    @(1943,2,0.107601058578867,0),
    @(2005.376850816456,0,0.0217387002767945,0.002329146458227983),
    @(2023.740959505003,0,0.03537527425185113,0.002855403302840001),
    @(2040.602289229372,0,0.00777583124659088,0.001710682874249994),
    @(2069.548786927769,0,0.02603643801600719,0.002698899062634892),
    @(2131,0,0.1007472841569682,0.00624789359112981),
    @(2139.901811344992,0,0.02837000420202299,0.004437766020428178),
    @(1934.30317654436,4,0.02278408802608079,0.002448225345049951),
    @(2183,0,0.02397241532204131,0),
    @(2183,2,0.02318956303588307,0),
    @(2201.968387101374,1,0.01957080465753616,0.002838310449775922),
    @(2207.71350700567,3,0.0169700664309912,0.003570376842864682),
    @(2232,1,0.05057059001395432,0.0009042413536712422),
    @(2273.578833093549,1,0.001930560522152485,0.0003210257614041017),
    @(2299.870800905546,1,0.004645426430644592,0.0003687385315539132),
    @(2413.274643654517,1,0.002700983249740335,0.0007528977333109513),
    @(2419.342384365287,1,0.008898478397206368,0.0006034180975013629),
    @(2432.48629598985,1,0.002089142602532552,0.0003043698656391923),
    @(2452.752171577677,1,0.02665842466937973,0.001125521300370366),
    @(2466.124342510111,3,0.005951846860492331,0.0005001594671462294),
    @(2478.166239524158,1,0.006306926592741646,0.0004048095903461701),
    @(2494.039958224013,0,0.0009677315022444048,0.001129020085951805),
    @(2514.257468531217,1,0.003029462879641493,0.0004095390246803904),
    @(2522.834832705841,1,0.0167442510284711,0.0009073335307082051),
    @(2532.957133048545,0,0.009265629304001955,0.002438323501053147),
    @(2544.746580038972,1,0.01030030361599892,0.0005973739577938043),
    @(2554.519429307579,1,0.006852808902499502,0.0005242383533031732),
)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

# Remove now-obsolete rows 49-53 (the table shrank from 53 to 48 data rows)
$ws.Range("A49:E53").Delete()
